$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column F (dSF) values for rows 2-28 per repulled data
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 10
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = 5
$ws.Range("F10").Value = 7
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = 4
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -1
$ws.Range("F23").Value = -6
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = -2
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 0
